$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''34.544.86'
$ws.Range("E2").Value = '  +1.19%  '
$ws.Range("D3").Value = '''1.796.21'
$ws.Range("E3").Value = '  +0.39%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''227.10'
$ws.Range("E5").Value = '  +0.33%  '
$ws.Range("D6").Value = '''0.557'
$ws.Range("E6").Value = '  +1.92%  '
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '''32.64'
$ws.Range("E8").Value = '  +2.20%  '
$ws.Range("E9").Value = '  +1.49%  '
$ws.Range("D10").Value = '''0.0695'
$ws.Range("E10").Value = '  +0.69%  '
$ws.Range("E11").Value = '  +0.35%  '
$ws.Range("D12").Value = '''2.055.52'
$ws.Range("E12").Value = '  +0.40%  '
$ws.Range("E13").Value = '  -0.64%  '
$ws.Range("D14").Value = '''1.799.61'
$ws.Range("E14").Value = '  +0.44%  '
$ws.Range("D15").Value = '''0.641'
$ws.Range("E15").Value = '  +3.31%  '
$ws.Range("D16").Value = '''34.482.35'
$ws.Range("E16").Value = '  +1.19%  '
$ws.Range("D17").Value = '''4.28'
$ws.Range("E17").Value = '  +2.29%  '
$ws.Range("D18").Value = '''68.83'
$ws.Range("E18").Value = '  +1.09%  '
$ws.Range("D19").Value = '''247.33'
$ws.Range("E19").Value = '  +0.94%  '
$ws.Range("D20").Value = '''0.0₃0802'
$ws.Range("E20").Value = '  +2.95%  '
$ws.Range("D21").Value = '''11.26'
$ws.Range("E21").Value = '  +3.99%  '
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("E23").Value = '  +1.89%  '
$ws.Range("D24").Value = '''2.08'
$ws.Range("E24").Value = '  +1.52%  '
$ws.Range("D25").Value = '''163.81'
$ws.Range("E25").Value = '  +1.47%  '
$ws.Range("D26").Value = '''7.27'
$ws.Range("E26").Value = '  +1.66%  '
$ws.Range("D27").Value = '''16.55'
$ws.Range("E27").Value = '  +1.53%  '
$ws.Range("E28").Value = '  +2.52%  '
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = '''1.24'
$ws.Range("E30").Value = '  +0.43%  '
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value = '''3.81'
$ws.Range("E31").Value = '  +3.87%  '
$ws.Range("D32").Value = '''0.0522'
$ws.Range("E32").Value = '  +0.79%  '
$ws.Range("D33").Value = '''3.89'
$ws.Range("E33").Value = '  +7.64%  '
$ws.Range("E34").Value = '  +1.43%  '
$ws.Range("D35").Value = '''1.444.93'
$ws.Range("E35").Value = '  -0.54%  '
$ws.Range("D36").Value = '''2.59'
$ws.Range("E36").Value = '  +6.91%  '
$ws.Range("E37").Value = '  +3.76%  '
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").Value = '''1.06'
$ws.Range("E38").Value = '  +2.28%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '''0.0192'
$ws.Range("E39").Value = '  +0.23%  '
$ws.Range("D40").Value = '''84.28'
$ws.Range("E40").Value = '  +4.99%  '
$ws.Range("E41").Value = '  +1.37%  '
$ws.Range("D42").Value = '''0.935'
$ws.Range("E42").Value = '  +1.59%  '
$ws.Range("E43").Value = '  +2.35%  '
$ws.Range("D44").Value = '''13.70'
$ws.Range("E44").Value = '  +1.48%  '
$ws.Range("E45").Value = '  +3.13%  '
$ws.Range("E46").Value = '  +0.84%  '
$ws.Range("E47").Value = '  +0.56%  '
$ws.Range("D48").Value = '''1.951.74'
$ws.Range("E48").Value = '  +0.17%  '
$ws.Range("D49").Value = '''105.94'
$ws.Range("E49").Value = '  -0.24%  '
$ws.Range("D50").Value = '''0.0₆0132'
$ws.Range("E50").Value = '  -1.94%  '
$ws.Range("E51").Value = '  +0.01%  '
